# Atualizacao de bases das ligas, do dia: 11-04-2024 as 23:56
#
# This script:
#  1) Swaps the full data payload (columns B..AC) between specific pairs of
#     adjacent rows (the "id" in column A stays put; every other column -
#     match id, div, date, teams, goals, result, odds... - moves with the row).
#  2) Performs a 3-way cyclic rotation of rows 221/222/223 (same B..AC payload).
#  3) Applies a handful of isolated odds corrections on rows 256-262.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-RowPayload {
    param($ws, $rowA, $rowB, $firstCol, $lastCol)

    $valsA = @{}
    $valsB = @{}
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $valsA[$c] = $ws.Cells.Item($rowA, $c).Value()
        $valsB[$c] = $ws.Cells.Item($rowB, $c).Value()
    }
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $ws.Cells.Item($rowA, $c).Value = $valsB[$c]
        $ws.Cells.Item($rowB, $c).Value = $valsA[$c]
    }
}

function Rotate-RowPayloadUp {
    # row[0] <- row[1] <- row[2] <- ... <- row[n-1] <- row[0]
    param($ws, $rows, $firstCol, $lastCol)

    $snapshot = @{}
    foreach ($r in $rows) {
        $rowVals = @{}
        for ($c = $firstCol; $c -le $lastCol; $c++) {
            $rowVals[$c] = $ws.Cells.Item($r, $c).Value()
        }
        $snapshot[$r] = $rowVals
    }

    $n = $rows.Length
    for ($i = 0; $i -lt $n; $i++) {
        $destRow = $rows[$i]
        $srcRow  = $rows[($i + 1) % $n]
        $srcVals = $snapshot[$srcRow]
        for ($c = $firstCol; $c -le $lastCol; $c++) {
            $ws.Cells.Item($destRow, $c).Value = $srcVals[$c]
        }
    }
}

# --- 1) pairwise swaps -------------------------------------------------
Swap-RowPayload $ws 74  75  2 29
Swap-RowPayload $ws 88  89  2 29
Swap-RowPayload $ws 111 112 2 29
Swap-RowPayload $ws 134 135 2 29
Swap-RowPayload $ws 140 141 2 29
Swap-RowPayload $ws 142 143 2 29
Swap-RowPayload $ws 167 168 2 29
Swap-RowPayload $ws 201 202 2 29

# --- 2) 3-way cyclic rotation (221 -> 222 -> 223 -> 221) ---------------
Rotate-RowPayloadUp $ws @(221, 222, 223) 2 29

# --- 3) isolated odds corrections --------------------------------------
$ws.Range("R256").Value = 1.925
$ws.Range("S256").Value = 1.925

$ws.Range("R257").Value = 2.05
$ws.Range("S257").Value = 1.8
$ws.Range("U257").Value = 2.05
$ws.Range("V257").Value = 1.8

$ws.Range("N258").Value = 2.3
$ws.Range("P258").Value = 3
$ws.Range("R258").Value = 2.025
$ws.Range("S258").Value = 1.825
$ws.Range("U258").Value = 2
$ws.Range("V258").Value = 1.85

$ws.Range("N260").Value = 2.15
$ws.Range("O260").Value = 3.1
$ws.Range("R260").Value = 1.875
$ws.Range("S260").Value = 1.975
$ws.Range("T260").Value = 2.25
$ws.Range("U260").Value = 1.95
$ws.Range("V260").Value = 1.9

$ws.Range("N261").Value = 1.6
$ws.Range("O261").Value = 3.8
$ws.Range("P261").Value = 5
$ws.Range("Q261").Value = -1
$ws.Range("R261").Value = 2.05
$ws.Range("S261").Value = 1.8

$ws.Range("R262").Value = 1.8
$ws.Range("S262").Value = 2.05
$ws.Range("U262").Value = 1.875
$ws.Range("V262").Value = 1.975
